$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'51.422.62"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -15.39%  '
$ws.Range('D3').Value = "'2.253.96"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -22.42%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'428.22"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -18.81%  '
$ws.Range('D6').Value = "'113.71"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -20.96%  '
$ws.Range('D7').Value = "'0.997"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').Value = "'0.447"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -18.13%  '
$ws.Range('D9').Value = "'2.264.01"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -22.25%  '
$ws.Range('D10').Value = "'5.11"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -15.63%  '
$ws.Range('D11').Value = "'0.0822"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -23.55%  '
$ws.Range('D12').Value = "'0.290"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -19.12%  '
$ws.Range('D13').Value = "'0.119"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -7.06%  '
$ws.Range('D14').Value = "'2.639.14"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -22.63%  '
$ws.Range('D15').Value = "'51.434.82"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -15.28%  '
$ws.Range('D16').Value = "'17.95"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -20.40%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = "'2.261.42"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -22.42%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = "'0.0000111"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -21.16%  '
$ws.Range('D19').Value = "'3.80"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -23.37%  '
$ws.Range('D20').Value = "'287.23"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -18.65%  '
$ws.Range('D21').Value = "'0.993"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.72%  '
$ws.Range('E22').Value = '  -0.77%  '
$ws.Range('D23').Value = "'8.31"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -28.15%  '
$ws.Range('D24').Value = "'4.83"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -26.05%  '
$ws.Range('D25').Value = "'0.997"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('D26').Value = "'52.08"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -19.85%  '
$ws.Range('D27').Value = "'0.353"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -21.82%  '
$ws.Range('D28').Value = "'2.333.20"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -23.05%  '
$ws.Range('D29').Value = "'0.132"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -25.20%  '
$ws.Range('D30').Value = "'0.996"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.29%  '
$ws.Range('D31').Value = "'6.55"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -16.61%  '
$ws.Range('D32').Value = "'142.62"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.66%  '
$ws.Range('D33').Value = "'0.0₃0607"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -29.14%  '
$ws.Range('D34').Value = "'16.10"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -17.75%  '
$ws.Range('D35').Value = "'1.27"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -24.10%  '
$ws.Range('D36').Value = "'4.50"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -19.26%  '
$ws.Range('D37').Value = "'0.998"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').Value = "'0.757"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -23.93%  '
$ws.Range('D39').Value = "'3.21"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -26.55%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = "'31.42"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -16.12%  '
$ws.Range('B41').Value = 'ImmutableX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D41').Value = "'0.940"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -21.28%  '
$ws.Range('D42').Value = "'10.11"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.28%  '
$ws.Range('D43').Value = "'0.546"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -16.00%  '
$ws.Range('D44').Value = "'0.0479"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -17.76%  '
$ws.Range('D45').Value = "'2.98"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -20.03%  '
$ws.Range('D46').Value = "'1.832.62"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -20.09%  '
$ws.Range('D47').Value = "'1.08"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -26.44%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = "'0.0779"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -14.90%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = "'0.0194"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -18.28%  '
$ws.Range('B50').Value = 'ZEEBU'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu'
$ws.Range('D50').Value = "'4.62"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.30%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = "'3.81"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -23.28%  '
